# Applies the odds updates described in the commit diff
# ("Atualizando o arquivo XLSX") to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 ---------------------------------------------------------------
$ws.Cells.Item(4, 7).Value  = 5        # G4  Odd_H_FT
$ws.Cells.Item(4, 8).Value  = 3.35     # H4  Odd_D_FT
$ws.Cells.Item(4, 9).Value  = 1.7      # I4  Odd_A_FT
$ws.Cells.Item(4, 10).Value = 5.1      # J4  Odd_H_HT
$ws.Cells.Item(4, 11).Value = 2.1      # K4  Odd_D_HT
$ws.Cells.Item(4, 12).Value = 2.25     # L4  Odd_A_HT

$ws.Cells.Item(4, 15).Value = 1.36     # O4  Odd_Over15_FT
$ws.Cells.Item(4, 16).Value = 2.7      # P4  Odd_Under15_FT
$ws.Cells.Item(4, 17).Value = 2.05     # Q4  Odd_Over25_FT
$ws.Cells.Item(4, 18).Value = 1.62     # R4  Odd_Under25_FT

$ws.Cells.Item(4, 21).Value = 1.98     # U4  Odd_BTTS_Yes
$ws.Cells.Item(4, 22).Value = 1.65     # V4  Odd_BTTS_No
$ws.Cells.Item(4, 23).Value = 11.75    # W4  Odd_CS_1-0
$ws.Cells.Item(4, 24).Value = 28       # X4  Odd_CS_2-0
$ws.Cells.Item(4, 25).Value = 16.5     # Y4  Odd_CS_2-1
$ws.Cells.Item(4, 26).Value = 100      # Z4  Odd_CS_3-0
$ws.Cells.Item(4, 27).Value = 55       # AA4 Odd_CS_3-1
$ws.Cells.Item(4, 28).Value = 60       # AB4 Odd_CS_3-2
$ws.Cells.Item(4, 29).Value = 8.25     # AC4 Odd_CS_0-0
$ws.Cells.Item(4, 30).Value = 6.7      # AD4 Odd_CS_1-1
$ws.Cells.Item(4, 31).Value = 18       # AE4 Odd_CS_2-2
$ws.Cells.Item(4, 32).Value = 100      # AF4 Odd_CS_3-3
$ws.Cells.Item(4, 33).Value = 1000     # AG4 Odd_CS_4-4
$ws.Cells.Item(4, 34).Value = 5.7      # AH4 Odd_CS_0-1
$ws.Cells.Item(4, 35).Value = 7.1      # AI4 Odd_CS_0-2
$ws.Cells.Item(4, 36).Value = 8.25     # AJ4 Odd_CS_1-2
$ws.Cells.Item(4, 37).Value = 12.5     # AK4 Odd_CS_0-3
$ws.Cells.Item(4, 38).Value = 15       # AL4 Odd_CS_1-3
$ws.Cells.Item(4, 39).Value = 32       # AM4 Odd_CS_2-3
$ws.Cells.Item(4, 40).Value = 6.5      # AN4 Odd_CS_1-0_HT
$ws.Cells.Item(4, 41).Value = 29       # AO4 Odd_CS_2-0_HT
$ws.Cells.Item(4, 42).Value = 35       # AP4 Odd_CS_2-1_HT
$ws.Cells.Item(4, 43).Value = 175      # AQ4 Odd_CS_3-0_HT
$ws.Cells.Item(4, 44).Value = 200      # AR4 Odd_CS_3-1_HT
$ws.Cells.Item(4, 45).Value = 450      # AS4 Odd_CS_3-2_HT

$ws.Cells.Item(4, 47).Value = 7.6      # AU4 Odd_CS_1-1_HT
$ws.Cells.Item(4, 48).Value = 75       # AV4 Odd_CS_2-2_HT
$ws.Cells.Item(4, 49).Value = 3.4      # AW4 Odd_CS_0-1_HT
$ws.Cells.Item(4, 50).Value = 8.25     # AX4 Odd_CS_0-2_HT
$ws.Cells.Item(4, 51).Value = 18.5     # AY4 Odd_CS_1-2_HT
$ws.Cells.Item(4, 52).Value = 28       # AZ4 Odd_CS_0-3_HT
$ws.Cells.Item(4, 53).Value = 65       # BA4 Odd_CS_1-3_HT
$ws.Cells.Item(4, 54).Value = 250      # BB4 Odd_CS_2-3_HT

# --- Row 6 ---------------------------------------------------------------
$ws.Cells.Item(6, 17).Value = 1.9      # Q6  Odd_Over25_FT
$ws.Cells.Item(6, 18).Value = 1.9      # R6  Odd_Under25_FT

# --- Row 9 ---------------------------------------------------------------
$ws.Cells.Item(9, 9).Value  = 3.3      # I9  Odd_A_FT
$ws.Cells.Item(9, 38).Value = 29       # AL9 Odd_CS_0-3
